# Wish List.xlsx - "2 more bugs added to wish list" edit
#
# The old rows 109-117 held a block of "new font system" follow-up bugs.
# The three that still had no replacement text (rows 109, 110, 111 - shared
# strings 165/166/167) are being repurposed for three brand-new bugs about
# replaying a level, and the now-stale six rows below them (112-117, shared
# strings 168-173: "Start New Game", "Select Profile State", "Delete Profile
# State", "Click ok state", "create profile state", "engame text") are being
# removed outright. A few of the spare blank rows further down are trimmed
# too, and row 110/111 lose their "Dave" / date tracking cells (111's note
# moves from column C to column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-purpose the three still-useful rows with the new bug text.
$ws.Range("C109").Value = "replay level breaks the instruction list highlighting for current level"
$ws.Range("C110").Value = "total score needs to be calculated on the fly based off of all your level scores"

# Row 111's note shifts from column C to column D, with new bug text.
$ws.Range("C111").Value = "replaying a level over and over keeps increasing total score"
$ws.Range("C111").Copy($ws.Range("D111"))
$ws.Range("C111").Clear()

# Row 110 no longer carries the "Dave" / date-assigned columns.
$ws.Range("A110:B110").Clear()

# Drop the six now-obsolete rows entirely (Start New Game / Select Profile
# State / Delete Profile State / Click ok state / create profile state /
# engame text).
$ws.Range("A112:A117").EntireRow.Delete()

# Trim three of the filler blank rows that followed them.
$ws.Range("A117:A119").EntireRow.Delete()

# Restore the view's selection to where it ends up after the edit.
$ws.Range("B113").Select()
